$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.836.29"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.534.94"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.78"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.91"
$ws.Range("E6").Value = "  +5.77%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.204"
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.651"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.73"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.54"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "4.094.04"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "603.99"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "69.987.53"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.73"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").Value = "3.539.19"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.19"
$ws.Range("E22").Value = "  +4.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.28"
$ws.Range("E23").Value = "  +5.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.41"
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  +5.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.63"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("E30").Value = "  +17.99%  "
$ws.Range("E31").Value = "  +2.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.57"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.24"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").Value = "0.0₃0860"
$ws.Range("E35").Value = "  +12.56%  "
$ws.Range("D36").Value = "3.746.09"
$ws.Range("E36").Value = "  +5.88%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.61"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.61"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "485.95"
$ws.Range("E42").Value = "  -7.51%  "
$ws.Range("E43").Value = "  -5.18%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("E45").Value = "  -2.48%  "
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.59"
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000254"
$ws.Range("E50").Value = "  +6.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.90"
$ws.Range("E51").Value = "  -0.77%  "
